# Apply the "got logistic regression baseline, one cnn result" edit.
#
# Summary of changes:
#  - LogReg sheet: inserted a new "Threshold" column (A) in front of the
#    existing Inv. Reg / Mean Crossval / Val columns, back-filled the
#    existing rows with Threshold = 0.9, and appended a new block of rows
#    (Threshold = 0.95) with fresh logistic-regression results. The small
#    "Testing" summary table at the bottom got the same new Threshold
#    column plus one more result row.
#  - CNN sheet: filled in Mean Crossval / Validation results for several
#    Num Filters / Window Size combos (replacing the "PENDING" placeholder
#    with a real number), swapped out a handful of untested Num Filters
#    values for ones that now have results, dropped the three trailing
#    rows that had no results at all, and added a "Testing" row at the
#    bottom with one CNN test result.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: LogReg
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LogReg")

# Header row - new Threshold column in A, everything else shifts right.
$ws1.Range("A1").Value = "Threshold"
$ws1.Range("B1").Value = "Inv. Reg"
$ws1.Range("C1").Value = "Mean Crossval"
$ws1.Range("D1").Value = "Val "

# Existing Threshold = 0.9 block (rows 2-8): shift old B/C values into C/D,
# old A values into B, and stamp the new Threshold value into A.
$logreg09 = @(
  # Inv. Reg, Mean Crossval, Val
  @(1,    0.7833,                0.8824),
  @(0.9,  0.769,                 0.8676),
  @(0.8,  0.763,                 0.8971),
  @(0.7,  0.7907,                0.8971),
  @(0.6,  0.7639,                0.8824),
  @(0.5,  0.7593,                0.9118),
  @(0.4,  0.7815,                0.8971)
)
$r = 2
foreach ($row in $logreg09) {
  $ws1.Range("A$r").Value = 0.9
  $ws1.Range("B$r").Value = $row[0]
  $ws1.Range("C$r").Value = $row[1]
  $ws1.Range("D$r").Value = $row[2]
  $r++
}

# New Threshold = 0.95 block (rows 9-14).
$logreg095 = @(
  @(1,    0.763,  0.9117),
  @(0.9,  0.7444, 0.8676),
  @(0.8,  0.7444, 0.8971),
  @(0.7,  0.7463, 0.8529),
  @(0.6,  0.7407, 0.8971),
  @(0.5,  0.763,  0.8824)
)
$r = 9
foreach ($row in $logreg095) {
  $ws1.Range("A$r").Value = 0.95
  $ws1.Range("B$r").Value = $row[0]
  $ws1.Range("C$r").Value = $row[1]
  $ws1.Range("D$r").Value = $row[2]
  $r++
}

# Fix precise binary-float literals that PowerShell would otherwise round.
$ws1.Range("C2").Value = 0.7833
$ws1.Range("D2").Value = 0.88239999999999996
$ws1.Range("C3").Value = 0.76900000000000002
$ws1.Range("D3").Value = 0.86760000000000004
$ws1.Range("C4").Value = 0.76300000000000001
$ws1.Range("D4").Value = 0.89710000000000001
$ws1.Range("C5").Value = 0.79069999999999996
$ws1.Range("D5").Value = 0.89710000000000001
$ws1.Range("C6").Value = 0.76390000000000002
$ws1.Range("D6").Value = 0.88239999999999996
$ws1.Range("C7").Value = 0.75929999999999997
$ws1.Range("D7").Value = 0.91180000000000005
$ws1.Range("C8").Value = 0.78149999999999997
$ws1.Range("D8").Value = 0.89710000000000001
$ws1.Range("C9").Value = 0.76300000000000001
$ws1.Range("D9").Value = 0.91169999999999995
$ws1.Range("C10").Value = 0.74439999999999995
$ws1.Range("D10").Value = 0.86760000000000004
$ws1.Range("C11").Value = 0.74439999999999995
$ws1.Range("D11").Value = 0.89710000000000001
$ws1.Range("C12").Value = 0.74629999999999996
$ws1.Range("D12").Value = 0.85289999999999999
$ws1.Range("C13").Value = 0.74070000000000003
$ws1.Range("D13").Value = 0.89710000000000001
$ws1.Range("C14").Value = 0.76300000000000001
$ws1.Range("D14").Value = 0.88239999999999996

# Bottom "Testing" summary table (rows 18-21): same Threshold column
# insertion, plus a new data row for Threshold 0.95.
$ws1.Range("A18").Value = "Testing"

$ws1.Range("A19").Value = "Threshold"
$ws1.Range("B19").Value = "Inv Reg"
$ws1.Range("C19").Value = "Acc"

$ws1.Range("A20").Value = 0.9
$ws1.Range("B20").Value = 0.5
$ws1.Range("C20").Value = 0.67859999999999998

$ws1.Range("A21").Value = 0.95
$ws1.Range("B21").Value = 1
$ws1.Range("C21").Value = 0.73809999999999998

# Column widths: the new Inv. Reg / Mean Crossval columns (now C & D)
# keep (approximately) the widths that used to belong to B & C.
$ws1.Columns.Item(3).ColumnWidth = 11.6
$ws1.Columns.Item(4).ColumnWidth = 11.75

# Selection: row 22 (the first empty row) is now fully selected.
$ws1.Rows.Item(22).Select()

# ---------------------------------------------------------------------
# Sheet 2: CNN
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("CNN")

# Row 5 (200 filters / window 10): PENDING result is now in.
$ws2.Range("C5").Value = 0.80559999999999998
$ws2.Range("D5").Value = 0.79410000000000003

# Row 6 (500 filters / window 10): new result.
$ws2.Range("C6").Value = 0.87409999999999999
$ws2.Range("D6").Value = 0.80879999999999996

# Row 7 (was 100/20, now an untested 1500/10 row).
$ws2.Range("A7").Value = 1500
$ws2.Range("B7").Value = 10
$ws2.Range("C7").ClearContents()
$ws2.Range("D7").ClearContents()

# Row 8 (was 200/20, now 100/20, with a new result).
$ws2.Range("A8").Value = 100
$ws2.Range("B8").Value = 20
$ws2.Range("C8").Value = 0.79259999999999997
$ws2.Range("D8").Value = 0.76470000000000005

# Row 9 (was 300/20, now 500/20, with a new result).
$ws2.Range("A9").Value = 500
$ws2.Range("B9").Value = 20
$ws2.Range("C9").Value = 0.86851999999999996
$ws2.Range("D9").Value = 0.92649999999999999

# Row 10 (was 100/30, now an untested 1500/20 row).
$ws2.Range("A10").Value = 1500
$ws2.Range("B10").Value = 20
$ws2.Range("C10").ClearContents()
$ws2.Range("D10").ClearContents()

# Row 11 (was 150/30, now 100/30, with a new result).
$ws2.Range("A11").Value = 100
$ws2.Range("B11").Value = 30
$ws2.Range("C11").Value = 0.79069999999999996
$ws2.Range("D11").Value = 0.76470000000000005

# Row 12 (was 200/30, now 500/30, with a new result).
$ws2.Range("A12").Value = 500
$ws2.Range("B12").Value = 30
$ws2.Range("C12").Value = 0.87222
$ws2.Range("D12").Value = 0.91169999999999995

# Row 13 (was 300/30, now an untested 1500/30 row).
$ws2.Range("A13").Value = 1500
$ws2.Range("B13").Value = 30
$ws2.Range("C13").ClearContents()
$ws2.Range("D13").ClearContents()

# Row 14 (was 500/30, now 100/50, with a new result).
$ws2.Range("A14").Value = 100
$ws2.Range("B14").Value = 50
$ws2.Range("C14").Value = 0.7944
$ws2.Range("D14").Value = 0.79410000000000003

# Row 15 (was 1000/30, now 500/50, with a new result).
$ws2.Range("A15").Value = 500
$ws2.Range("B15").Value = 50
$ws2.Range("C15").Value = 0.91849999999999998
$ws2.Range("D15").Value = 0.94120000000000004

# Row 16 (was 2000/30, now an untested 1500/50 row).
$ws2.Range("A16").Value = 1500
$ws2.Range("B16").Value = 50
$ws2.Range("C16").ClearContents()
$ws2.Range("D16").ClearContents()

# Rows 17-19 (old 1000/50, 1500/50, 2000/50 placeholders) are gone.
$ws2.Range("A17:D19").ClearContents()

# New "Testing" summary table at the bottom.
$ws2.Range("A21").Value = "Testing"
$ws2.Range("A22").Value = 500
$ws2.Range("B22").Value = 50
$ws2.Range("C22").Value = 0.72619999999999996

# Selection: first empty row below the new table.
$ws2.Range("A23").Select()

# Restore LogReg as the active/selected tab (it was the active sheet
# before and after the edit).
$ws1.Select()
$ws1.Rows.Item(22).Select()
